$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the old row 13, shifting the
# existing rows 13-29 down to 14-30 (dimension grows from A1:R29 to A1:R30).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new week's data.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44721
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("N13").Value = "$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 833
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"
